$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 2998.75  # was 3498.3333
$ws.Range("I20").Value = 2998.75  # was 3498.3333
$ws.Range("K20").Value = 2998.75  # was 3498.3333
$ws.Range("M20").Value = -2768.75  # was -3268.3333
$ws.Range("H35").Value = 2998.75  # was 3498.3333
$ws.Range("I35").Value = 2998.75  # was 3498.3333
$ws.Range("K35").Value = 2998.75  # was 3498.3333
$ws.Range("M35").Value = -2619.75  # was -3119.3333
$ws.Range("H106").Value = 11992  # was 9324.666999999999
$ws.Range("I106").Value = 11992  # was 9324.666999999999
$ws.Range("K106").Value = 11992  # was 9324.666999999999
$ws.Range("M106").Value = -11361  # was -8693.666999999999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1060221.8  # was 1060228.4
$ws.Range("I61").Value = 3724.0527  # was 3734.5264
$ws.Range("K61").Value = 3724.0527  # was 3734.5264
$ws.Range("M61").Value = -3512.0527  # was -3522.5264
$ws.Range("H74").Value = 22141  # was 20094.094
$ws.Range("I74").Value = 1882.3684  # was 1667.591
$ws.Range("K74").Value = 1882.3684  # was 1667.591
$ws.Range("M74").Value = -1008.3684  # was -793.5909999999999
$ws.Range("H77").Value = 22141  # was 20094.094
$ws.Range("I77").Value = 1882.3684  # was 1667.591
$ws.Range("K77").Value = 9411.842000000001  # was 8337.955
$ws.Range("M77").Value = -5043.842000000001  # was -3969.955
$ws.Range("H136").Value = 1060221.8  # was 1060228.4
$ws.Range("I136").Value = 3724.0527  # was 3734.5264
$ws.Range("K136").Value = 11172.1581  # was 11203.5792
$ws.Range("M136").Value = -8622.158100000001  # was -8653.5792

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H23").Value = 8259.75  # was 0
$ws.Range("J23").Value = 8259.75  # was 0
$ws.Range("L23").Value = 8259.75  # was 0
$ws.Range("N23").Value = -8825.75  # new cell
$ws.Range("H31").Value = 14975  # was 7499
$ws.Range("I31").Value = 0  # was 23
$ws.Range("K31").Value = 0  # was 23
$ws.Range("M31").ClearContents()  # was 229
$ws.Range("H36").Value = 765.3333  # was 780.5
$ws.Range("I36").Value = 765.3333  # was 780.5
$ws.Range("K36").Value = 765.3333  # was 780.5
$ws.Range("M36").Value = -231.3333  # was -246.5
$ws.Range("H38").Value = 5833  # was 0
$ws.Range("J38").Value = 5833  # was 0
$ws.Range("L38").Value = 5833  # was 0
$ws.Range("N38").Value = -6665  # new cell
$ws.Range("H39").Value = 20000  # was 3200
$ws.Range("I39").Value = 0  # was 3200
$ws.Range("J39").Value = 20000  # was 0
$ws.Range("K39").Value = 0  # was 3200
$ws.Range("L39").Value = 20000  # was 0
$ws.Range("M39").ClearContents()  # was -2811
$ws.Range("N39").Value = -20778  # new cell
$ws.Range("H99").Value = 14139.895  # was 12983.619
$ws.Range("I99").Value = 15803.625  # was 14973.941
$ws.Range("J99").Value = 5266.6665  # was 4524.75
$ws.Range("K99").Value = 15803.625  # was 14973.941
$ws.Range("L99").Value = 5266.6665  # was 4524.75
$ws.Range("M99").Value = -14305.625  # was -13475.941
$ws.Range("N99").Value = -8262.666499999999  # was -7520.75
$ws.Range("H107").Value = 1927.5883  # was 1945.1765
$ws.Range("I107").Value = 2101.9  # was 2131.8
$ws.Range("K107").Value = 2101.9  # was 2131.8
$ws.Range("M107").Value = -181.9000000000001  # was -211.8000000000002
$ws.Range("H134").Value = 67760  # was 81019.39999999999
$ws.Range("I134").Value = 72331.28999999999  # was 84270.75
$ws.Range("J134").Value = 51760.5  # was 68014
$ws.Range("K134").Value = 216993.87  # was 252812.25
$ws.Range("L134").Value = 155281.5  # was 204042
$ws.Range("M134").Value = -214458.87  # was -250277.25
$ws.Range("N134").Value = -160351.5  # was -209112

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2954.8572  # was 2948.25
$ws.Range("I99").Value = 1686.9  # was 1708.8889
$ws.Range("J99").Value = 6124.75  # was 6666.3335
$ws.Range("K99").Value = 1686.9  # was 1708.8889
$ws.Range("L99").Value = 6124.75  # was 6666.3335
$ws.Range("M99").Value = -188.9000000000001  # was -210.8888999999999
$ws.Range("N99").Value = -9120.75  # was -9662.333500000001
$ws.Range("H110").Value = 250000  # was 121666.664
$ws.Range("J110").Value = 250000  # was 121666.664
$ws.Range("L110").Value = 250000  # was 121666.664
$ws.Range("N110").Value = -258180  # was -129846.664
$ws.Range("H126").Value = 2954.8572  # was 2948.25
$ws.Range("I126").Value = 1686.9  # was 1708.8889
$ws.Range("J126").Value = 6124.75  # was 6666.3335
$ws.Range("K126").Value = 5060.700000000001  # was 5126.6667
$ws.Range("L126").Value = 18374.25  # was 19999.0005
$ws.Range("M126").Value = -2590.700000000001  # was -2656.6667
$ws.Range("N126").Value = -23314.25  # was -24939.0005

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 491957.38  # was 91129.2
$ws.Range("I9").Value = 772662.1  # was 102023.5
$ws.Range("J9").Value = 724  # was 83866.336
$ws.Range("K9").Value = 2317986.3  # was 306070.5
$ws.Range("L9").Value = 2172  # was 251599.008
$ws.Range("M9").Value = -2317762.3  # was -305846.5
$ws.Range("N9").Value = -2620  # was -252047.008
$ws.Range("H68").Value = 869.25  # was 956.125
$ws.Range("J68").Value = 1029  # was 1144.8334
$ws.Range("L68").Value = 3087  # was 3434.5002
$ws.Range("N68").Value = -4709  # was -5056.5002
$ws.Range("H71").Value = 869.25  # was 956.125
$ws.Range("J71").Value = 1029  # was 1144.8334
$ws.Range("L71").Value = 9261  # was 10303.5006
$ws.Range("N71").Value = -17373  # was -18415.5006
$ws.Range("H86").Value = 542.6818  # was 556.6667
$ws.Range("J86").Value = 411.66666  # was 426.45456
$ws.Range("L86").Value = 1234.99998  # was 1279.36368
$ws.Range("N86").Value = -3606.99998  # was -3651.36368
$ws.Range("H89").Value = 542.6818  # was 556.6667
$ws.Range("J89").Value = 411.66666  # was 426.45456
$ws.Range("L89").Value = 3704.99994  # was 3838.09104
$ws.Range("N89").Value = -15560.99994  # was -15694.09104
$ws.Range("H92").Value = 299.5  # was 304
$ws.Range("J92").Value = 292.25  # was 303.5
$ws.Range("L92").Value = 876.75  # was 910.5
$ws.Range("N92").Value = -3372.75  # was -3406.5
$ws.Range("H122").Value = 9359303  # was 10250638
$ws.Range("J122").Value = 1577456.2  # was 1774601.5
$ws.Range("L122").Value = 14197105.8  # was 15971413.5
$ws.Range("N122").Value = -14202005.8  # was -15976313.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 504068.3  # was 489697.84
$ws.Range("I132").Value = 2771  # was 2717.1292
$ws.Range("K132").Value = 8313  # was 8151.3876
$ws.Range("M132").Value = -5783  # was -5621.3876

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2969.3  # was 2595.25
$ws.Range("I46").Value = 500  # was 650
$ws.Range("K46").Value = 500  # was 650
$ws.Range("M46").Value = -312  # was -462
$ws.Range("H93").Value = 4376.0356  # was 3882.7778
$ws.Range("I93").Value = 2797.6667  # was 3201.9583
$ws.Range("J93").Value = 13846.25  # was 9329.333000000001
$ws.Range("K93").Value = 2797.6667  # was 3201.9583
$ws.Range("L93").Value = 13846.25  # was 9329.333000000001
$ws.Range("M93").Value = -1549.6667  # was -1953.9583
$ws.Range("N93").Value = -16342.25  # was -11825.333
$ws.Range("H100").Value = 3488.9285  # was 2949.889
$ws.Range("I100").Value = 2956.25  # was 2160.4
$ws.Range("J100").Value = 4199.1665  # was 3936.75
$ws.Range("K100").Value = 2956.25  # was 2160.4
$ws.Range("L100").Value = 4199.1665  # was 3936.75
$ws.Range("M100").Value = -2415.25  # was -1619.4
$ws.Range("N100").Value = -5281.1665  # was -5018.75

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 40249.668  # was 48102.6
$ws.Range("I69").Value = 0  # was 60246
$ws.Range("J69").Value = 40249.668  # was 40007
$ws.Range("K69").Value = 0  # was 60246
$ws.Range("L69").Value = 40249.668  # was 40007
$ws.Range("M69").ClearContents()  # was -59497
$ws.Range("N69").Value = -41747.668  # was -41505
$ws.Range("H72").Value = 40249.668  # was 48102.6
$ws.Range("I72").Value = 0  # was 60246
$ws.Range("J72").Value = 40249.668  # was 40007
$ws.Range("K72").Value = 0  # was 180738
$ws.Range("L72").Value = 120749.004  # was 120021
$ws.Range("M72").ClearContents()  # was -176994
$ws.Range("N72").Value = -128237.004  # was -127509
$ws.Range("H103").Value = 22912  # was 36087.75
$ws.Range("J103").Value = 22912  # was 36087.75
$ws.Range("L103").Value = 22912  # was 36087.75
$ws.Range("N103").Value = -25256  # was -38431.75
$ws.Range("H111").Value = 0  # was 52750
$ws.Range("J111").Value = 0  # was 52750
$ws.Range("L111").Value = 0  # was 52750
$ws.Range("N111").ClearContents()  # was -60930
$ws.Range("H132").Value = 367118.66  # was 379778.6
$ws.Range("I132").Value = 1646.44  # was 1715.8334
$ws.Range("K132").Value = 4939.32  # was 5147.5002
$ws.Range("M132").Value = -2409.32  # was -2617.5002
